$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 0.630749404103418
$ws.Range("J2").Value = 0.630749404103418
$ws.Range("M2").Value = 0.8077876666666667
$ws.Range("N2").Value = 2.423363
$ws.Range("O2").Value = 0.04902038147436601
$ws.Range("P2").Value = 0.04902038147436601
$ws.Range("Q2").Value = 0.07445917448777778
$ws.Range("R2").Value = 0.6701325703900001
$ws.Range("S2").Value = 0.03091957640387859
$ws.Range("T2").Value = 0.03091957640387859
$ws.Range("I3").Value = 0.630749404103418
$ws.Range("J3").Value = 0.630749404103418
$ws.Range("O3").Value = 0.4722854529078861
$ws.Range("P3").Value = 0.4722854529078861
$ws.Range("S3").Value = 0.297893767988362
$ws.Range("T3").Value = 0.297893767988362
$ws.Range("I4").Value = 0.630749404103418
$ws.Range("J4").Value = 0.630749404103418
$ws.Range("M4").Value = 5.009781333333333
$ws.Range("N4").Value = 15.029344
$ws.Range("O4").Value = 0.3040172587389813
$ws.Range("P4").Value = 0.3040172587389813
$ws.Range("Q4").Value = 0.4617849440355555
$ws.Range("R4").Value = 4.156064496319999
$ws.Range("S4").Value = 0.1917587047867671
$ws.Range("T4").Value = 0.1917587047867671
$ws.Range("I5").Value = 0.630749404103418
$ws.Range("J5").Value = 0.630749404103418
$ws.Range("M5").Value = 2.878432333333333
$ws.Range("N5").Value = 8.635297
$ws.Range("O5").Value = 0.1746769068787666
$ws.Range("P5").Value = 0.1746769068787666
$ws.Range("Q5").Value = 0.2653242977122222
$ws.Range("R5").Value = 2.38791867941
$ws.Range("S5").Value = 0.1101773549244102
$ws.Range("T5").Value = 0.1101773549244102
$ws.Range("G6").Value = 0.05396166666666666
$ws.Range("H6").Value = 0.161885
$ws.Range("I6").Value = 0.369250595896582
$ws.Range("J6").Value = 0.369250595896582
$ws.Range("M6").Value = 0.8077876666666667
$ws.Range("N6").Value = 2.423363
$ws.Range("O6").Value = 0.04902038147436601
$ws.Range("P6").Value = 0.04902038147436601
$ws.Range("Q6").Value = 0.04358956880611111
$ws.Range("R6").Value = 0.392306119255
$ws.Range("S6").Value = 0.01810080507048742
$ws.Range("T6").Value = 0.01810080507048742
$ws.Range("G7").Value = 0.05396166666666666
$ws.Range("H7").Value = 0.161885
$ws.Range("I7").Value = 0.369250595896582
$ws.Range("J7").Value = 0.369250595896582
$ws.Range("O7").Value = 0.4722854529078861
$ws.Range("P7").Value = 0.4722854529078861
$ws.Range("Q7").Value = 0.4199624447316666
$ws.Range("R7").Value = 3.779662002585
$ws.Range("S7").Value = 0.1743916849195241
$ws.Range("T7").Value = 0.1743916849195241
$ws.Range("G8").Value = 0.05396166666666666
$ws.Range("H8").Value = 0.161885
$ws.Range("I8").Value = 0.369250595896582
$ws.Range("J8").Value = 0.369250595896582
$ws.Range("M8").Value = 5.009781333333333
$ws.Range("N8").Value = 15.029344
$ws.Range("O8").Value = 0.3040172587389813
$ws.Range("P8").Value = 0.3040172587389813
$ws.Range("Q8").Value = 0.2703361503822222
$ws.Range("R8").Value = 2.43302535344
$ws.Range("S8").Value = 0.1122585539522142
$ws.Range("T8").Value = 0.1122585539522142
$ws.Range("G9").Value = 0.05396166666666666
$ws.Range("H9").Value = 0.161885
$ws.Range("I9").Value = 0.369250595896582
$ws.Range("J9").Value = 0.369250595896582
$ws.Range("M9").Value = 2.878432333333333
$ws.Range("N9").Value = 8.635297
$ws.Range("O9").Value = 0.1746769068787666
$ws.Range("P9").Value = 0.1746769068787666
$ws.Range("Q9").Value = 0.1553250060938889
$ws.Range("R9").Value = 1.397925054845
$ws.Range("S9").Value = 0.0644995519543563
$ws.Range("T9").Value = 0.0644995519543563
